$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "26.632.44"
$ws.Cells.Item(2,5).Value = "  -0.26%  "
$ws.Cells.Item(3,4).Value = "1.597.12"
$ws.Cells.Item(3,5).Value = "  -0.20%  "
$ws.Cells.Item(4,4).Formula = "'1.00"
$ws.Cells.Item(4,5).Value = "  -0.03%  "
$ws.Cells.Item(5,4).Formula = "'210.59"
$ws.Cells.Item(5,5).Value = "  -0.37%  "
$ws.Cells.Item(6,5).Value = "  -0.52%  "
$ws.Cells.Item(7,4).Formula = "'1.00"
$ws.Cells.Item(7,5).Value = "  +0.01%  "
$ws.Cells.Item(8,5).Value = "  -0.70%  "
$ws.Cells.Item(9,5).Value = "  -0.37%  "
$ws.Cells.Item(10,4).Formula = "'19.56"
$ws.Cells.Item(10,5).Value = "  +0.01%  "
$ws.Cells.Item(11,4).Formula = "'0.0844"
$ws.Cells.Item(11,5).Value = "  +0.06%  "
$ws.Cells.Item(12,4).Value = "1.821.42"
$ws.Cells.Item(12,5).Value = "  -0.19%  "
$ws.Cells.Item(13,4).Value = "1.596.72"
$ws.Cells.Item(13,5).Value = "  -1.55%  "
$ws.Cells.Item(14,5).Value = "  -0.27%  "
$ws.Cells.Item(15,5).Value = "  -0.09%  "
$ws.Cells.Item(16,5).Value = "  -1.17%  "
$ws.Cells.Item(17,4).Value = "26.599.73"
$ws.Cells.Item(17,5).Value = "  -0.31%  "
$ws.Cells.Item(18,4).Value = "0.0₃0739"
$ws.Cells.Item(18,5).Value = "  -2.79%  "
$ws.Cells.Item(19,4).Formula = "'1.00"
$ws.Cells.Item(19,5).Value = "  +0.00%  "
$ws.Cells.Item(20,4).Formula = "'208.28"
$ws.Cells.Item(20,5).Value = "  -0.56%  "
$ws.Cells.Item(21,5).Value = "  -1.35%  "
$ws.Cells.Item(22,4).Formula = "'4.28"
$ws.Cells.Item(22,5).Value = "  +0.03%  "
$ws.Cells.Item(23,5).Value = "  -3.29%  "
$ws.Cells.Item(24,4).Formula = "'8.96"
$ws.Cells.Item(24,5).Value = "  +0.04%  "
$ws.Cells.Item(25,4).Formula = "'143.80"
$ws.Cells.Item(25,5).Value = "  +0.68%  "
$ws.Cells.Item(26,5).Value = "  -0.02%  "
$ws.Cells.Item(27,4).Formula = "'7.13"
$ws.Cells.Item(27,5).Value = "  +0.15%  "
$ws.Cells.Item(28,4).Formula = "'0.113"
$ws.Cells.Item(28,5).Value = "  -1.05%  "
$ws.Cells.Item(29,4).Formula = "'15.26"
$ws.Cells.Item(29,5).Value = "  -0.66%  "
$ws.Cells.Item(30,5).Value = "  -2.37%  "
$ws.Cells.Item(31,5).Value = "  -0.59%  "
$ws.Cells.Item(32,4).Formula = "'3.25"
$ws.Cells.Item(32,5).Value = "  +0.00%  "
$ws.Cells.Item(33,5).Value = "  -0.40%  "
$ws.Cells.Item(34,4).Value = "1.278.70"
$ws.Cells.Item(34,5).Value = "  -1.12%  "
$ws.Cells.Item(35,2).Value = "HuobiToken"
$ws.Cells.Item(35,3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(35,4).Formula = "'2.48"
$ws.Cells.Item(35,5).Value = "  +0.53%  "
$ws.Cells.Item(36,2).Value = "WEMIXToken"
$ws.Cells.Item(36,3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(36,4).Formula = "'1.23"
$ws.Cells.Item(36,5).Value = "  +14.06%  "
$ws.Cells.Item(37,4).Formula = "'0.600"
$ws.Cells.Item(37,5).Value = "  -3.34%  "
$ws.Cells.Item(38,5).Value = "  -1.02%  "
$ws.Cells.Item(39,5).Value = "  -2.11%  "
$ws.Cells.Item(40,5).Value = "  -0.33%  "
$ws.Cells.Item(41,5).Value = "  +0.11%  "
$ws.Cells.Item(42,5).Value = "  -1.26%  "
$ws.Cells.Item(43,4).Formula = "'0.770"
$ws.Cells.Item(43,5).Value = "  -1.89%  "
$ws.Cells.Item(44,4).Formula = "'62.76"
$ws.Cells.Item(44,5).Value = "  -0.69%  "
$ws.Cells.Item(45,4).Value = "1.732.95"
$ws.Cells.Item(45,5).Value = "  -0.23%  "
$ws.Cells.Item(46,4).Formula = "'89.45"
$ws.Cells.Item(46,5).Value = "  -1.79%  "
$ws.Cells.Item(47,5).Value = "  -0.04%  "
$ws.Cells.Item(48,5).Value = "  +2.13%  "
$ws.Cells.Item(49,5).Value = "  +0.49%  "
$ws.Cells.Item(50,2).Value = "EnergySwap"
$ws.Cells.Item(50,3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(50,4).Formula = "'7.48"
$ws.Cells.Item(50,5).Value = "  +1.05%  "
$ws.Cells.Item(51,2).Value = "USDD"
$ws.Cells.Item(51,3).Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Cells.Item(51,4).Formula = "'1.00"
$ws.Cells.Item(51,5).Value = "  +0.01%  "
